$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=18; I='aa'; J='Agree/Accept'}
    @{Row=27; I='sv'; J='Statement-opinion'}
    @{Row=36; I='sd'; J='Statement-non-opinion'}
    @{Row=41; I='sd'; J='Statement-non-opinion'}
    @{Row=52; I='%'; J='Uninterpretable'}
    @{Row=74; I='sd'; J='Statement-non-opinion'}
    @{Row=78; I='%'; J='Uninterpretable'}
    @{Row=79; I='ba'; J='Appreciation'}
    @{Row=85; I='sd'; J='Statement-non-opinion'}
    @{Row=86; I='sd'; J='Statement-non-opinion'}
    @{Row=91; I='sd'; J='Statement-non-opinion'}
    @{Row=92; I='%'; J='Uninterpretable'}
    @{Row=96; I='sd'; J='Statement-non-opinion'}
    @{Row=114; I='sd'; J='Statement-non-opinion'}
    @{Row=124; I='sv'; J='Statement-opinion'}
    @{Row=130; I='sd'; J='Statement-non-opinion'}
    @{Row=150; I='sv'; J='Statement-opinion'}
    @{Row=153; I='aa'; J='Agree/Accept'}
    @{Row=160; I='sd'; J='Statement-non-opinion'}
    @{Row=166; I='aa'; J='Agree/Accept'}
    @{Row=167; I='aa'; J='Agree/Accept'}
    @{Row=168; I='aa'; J='Agree/Accept'}
    @{Row=169; I='sd'; J='Statement-non-opinion'}
    @{Row=174; I='aa'; J='Agree/Accept'}
    @{Row=175; I='aa'; J='Agree/Accept'}
    @{Row=184; I='%'; J='Uninterpretable'}
    @{Row=186; I='%'; J='Uninterpretable'}
    @{Row=187; I='sd'; J='Statement-non-opinion'}
    @{Row=189; I='%'; J='Uninterpretable'}
    @{Row=190; I='sd'; J='Statement-non-opinion'}
    @{Row=195; I='sd'; J='Statement-non-opinion'}
    @{Row=204; I='ba'; J='Appreciation'}
    @{Row=213; I='sd'; J='Statement-non-opinion'}
    @{Row=216; I='sd'; J='Statement-non-opinion'}
    @{Row=217; I='sd'; J='Statement-non-opinion'}
    @{Row=228; I='sd'; J='Statement-non-opinion'}
    @{Row=234; I='sd'; J='Statement-non-opinion'}
    @{Row=240; I='sv'; J='Statement-opinion'}
    @{Row=255; I='sd'; J='Statement-non-opinion'}
    @{Row=280; I='%'; J='Uninterpretable'}
    @{Row=292; I='ba'; J='Appreciation'}
    @{Row=301; I='sv'; J='Statement-opinion'}
    @{Row=313; I='sd'; J='Statement-non-opinion'}
    @{Row=314; I='aa'; J='Agree/Accept'}
    @{Row=338; I='sv'; J='Statement-opinion'}
    @{Row=358; I='sv'; J='Statement-opinion'}
    @{Row=368; I='sd'; J='Statement-non-opinion'}
    @{Row=371; I='ba'; J='Appreciation'}
    @{Row=383; I='sd'; J='Statement-non-opinion'}
    @{Row=386; I='sv'; J='Statement-opinion'}
    @{Row=418; I='%'; J='Uninterpretable'}
    @{Row=439; I='sd'; J='Statement-non-opinion'}
    @{Row=441; I='sd'; J='Statement-non-opinion'}
    @{Row=448; I='aa'; J='Agree/Accept'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows"